$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-03 Sunday" "2024-03-04 Monday"

Replace-Text "501÷4=" "660÷5="
Replace-Text "616÷8=" "825÷4="
Replace-Text "687÷8=" "953÷4="
Replace-Text "622÷8=" "502÷9="
Replace-Text "543÷2=" "854÷7="

Replace-Text "306÷4=" "297÷3="
Replace-Text "629÷3=" "837÷5="
Replace-Text "974÷6=" "225÷5="
Replace-Text "316÷5=" "351÷9="
Replace-Text "653÷4=" "978÷2="

Replace-Text "610÷6=" "657÷2="
Replace-Text "362÷9=" "366÷6="
Replace-Text "796÷8=" "295÷2="
Replace-Text "272÷6=" "843÷3="
Replace-Text "205÷5=" "803÷5="

Replace-Text "615÷5=" "638÷6="
Replace-Text "353÷4=" "177÷6="
Replace-Text "757÷6=" "374÷4="
Replace-Text "971÷7=" "168÷5="
Replace-Text "410÷2=" "143÷8="

Replace-Text "866÷5=" "664÷8="
Replace-Text "794÷7=" "105÷7="
Replace-Text "163÷7=" "100÷9="
Replace-Text "831÷4=" "282÷2="
Replace-Text "541÷8=" "239÷2="
